$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: RegNum KR66 NUO, Make SUZUKI, Colour BLUE->RED, Result Pass->Fail
$ws.Range("C7").Value = "RED"
$ws.Range("D7").Value = "Fail"

# Update the selection to match the diff (active cell C8)
$ws.Range("C8").Select()
